# Auto-generated edit script for cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Price" (D) column cells to Text format before writing so that
# values such as "137.01" or "0.999" are not auto-converted to numbers by Excel,
# matching the original inline-string cell type used in the workbook.
$dCells = @()
$dCells += "D2"
$dCells += "D3"
$dCells += "D5"
$dCells += "D6"
$dCells += "D9"
$dCells += "D12"
$dCells += "D14"
$dCells += "D15"
$dCells += "D17"
$dCells += "D18"
$dCells += "D19"
$dCells += "D20"
$dCells += "D21"
$dCells += "D24"
$dCells += "D25"
$dCells += "D26"
$dCells += "D27"
$dCells += "D29"
$dCells += "D30"
$dCells += "D31"
$dCells += "D32"
$dCells += "D34"
$dCells += "D36"
$dCells += "D37"
$dCells += "D43"
$dCells += "D46"

foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Set the new Price values (column D)
$ws.Range("D2").Value = "60.761.41"
$ws.Range("D3").Value = "2.371.69"
$ws.Range("D5").Value = "544.26"
$ws.Range("D6").Value = "141.16"
$ws.Range("D9").Value = "2.371.95"
$ws.Range("D12").Value = "5.32"
$ws.Range("D14").Value = "25.51"
$ws.Range("D15").Value = "2.800.02"
$ws.Range("D17").Value = "60.542.09"
$ws.Range("D18").Value = "2.374.14"
$ws.Range("D19").Value = "10.67"
$ws.Range("D20").Value = "4.11"
$ws.Range("D21").Value = "316.82"
$ws.Range("D24").Value = "1.84"
$ws.Range("D25").Value = "62.97"
$ws.Range("D26").Value = "0.999"
$ws.Range("D27").Value = "2.491.17"
$ws.Range("D29").Value = "7.81"
$ws.Range("D30").Value = "520.55"
$ws.Range("D31").Value = "1.43"
$ws.Range("D32").Value = "8.01"
$ws.Range("D34").Value = "1.84"
$ws.Range("D36").Value = "0.999"
$ws.Range("D37").Value = "4.65"
$ws.Range("D43").Value = "137.01"
$ws.Range("D46").Value = "139.80"

# Restore the original (default/Normal) cell style on the Price cells so only
# the text content differs, not the formatting, matching the source diff.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Set the new Volume(1h) values (column E) - these already contain padding
# spaces/percent signs so Excel keeps them as text automatically.
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -8.61%  "
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +4.19%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -3.19%  "
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("E29").Value = "  +3.60%  "
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("E38").Value = "  -5.70%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -5.39%  "
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -2.58%  "
